$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A7").Value = "value to be stash"
$ws.Range("A7").Select() | Out-Null
